# Revise LHS sample to uniform distribution, include design params.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (he_av_design): add Category "p" and Signs 1
$ws.Range("F8").Value = "p"
$ws.Range("G8").Value = 1

# Row 12 (rec_fr): Type -> 0, Distribution "pert" -> "uniform"
$ws.Range("E12").Value = 0
$ws.Range("H12").Value = "uniform"

# Row 13 (ab_rec): add Category "p" and Signs -1
$ws.Range("F13").Value = "p"
$ws.Range("G13").Value = -1

# Row 20 (eff_blk): Type -> 0, clear Distribution "normal" and B1 value
$ws.Range("E20").Value = 0
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = ""

# Reselect cell H13, matching the saved selection state in the sheet.
$ws.Range("H13").Select() | Out-Null
